$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 175, shifting existing rows 175..299 down to 176..300
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new data point
$ws.Range("A175").Value = 10
$ws.Range("B175").Value = "Vega Modelo de Temuco"
$ws.Range("C175").Value = "La Araucanía"
$ws.Range("D175").Value = 44907
$ws.Range("E175").Value = 9
$ws.Range("F175").Value = 100112052
$ws.Range("G175").Value = "Albahaca"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 80
$ws.Range("K175").Value = 8000
$ws.Range("L175").Value = 8000
$ws.Range("M175").Value = 8000
$ws.Range("N175").Value = "$/paquete"
$ws.Range("O175").Value = "Región del Maule"
$ws.Range("P175").Value = 8000
$ws.Range("Q175").Value = 1
$ws.Range("R175").Value = "Hortaliza"
